$wb = $excel.ActiveWorkbook

$capSheet = $wb.Worksheets.Item("Capacità di trasmissione MW")

# Duplicate "Capacità di trasmissione MW" and drop the copy right before it
# (i.e. as the first sheet in the workbook). Copying - rather than adding a
# blank sheet - means the new sheet inherits identical formatting, most
# notably the bestFit width already set on column A (which will also hold
# region names on the new sheet).
$capSheet.Copy($capSheet)
$newSheet = $wb.Worksheets.Item(1)
$newSheet.Name = "Info geografiche"

# Strip out everything the copy brought along except column A (the region
# names), including the old TRANSPOSE array formula that used to live in
# B1:H1.
$newSheet.Range("B1:H8").Clear() | Out-Null

# B1 was the anchor of that old spilled array formula; writing into it once
# fully detaches the leftover spill metadata so cell C1 can be written
# normally afterwards.
$newSheet.Range("B1").Value = "reset"

# Fill in the header row. Write C1 before B1 so that "latitude" claims a
# lower shared-string index than "longitude", matching source order.
$newSheet.Range("C1").Value = "latitude"
$newSheet.Range("B1").Value = "longitude"
$newSheet.Range("D1").Value = "altitude"

# Re-fetch the "Capacità di trasmissione MW" sheet by name: the reference
# obtained before the Copy()/rename above can end up stale.
$capSheet = $wb.Worksheets.Item("Capacità di trasmissione MW")
$capSheet.Activate() | Out-Null
$capSheet.Range("A2:A8").Select() | Out-Null

# Leave the new "Info geografiche" sheet selected/active, with E3 selected.
$newSheet.Activate() | Out-Null
$newSheet.Range("E3").Select() | Out-Null
